# Apply the changes described by the diff:
# 1. About sheet: C1 date value 45320 -> 45392
# 2. MCF sheet: update several Maximum Capacity Factor values to 1
# 3. Update the active selection on the MCF sheet to B17

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsMCF = $wb.Worksheets.Item("MCF")

# --- About sheet: update the "last updated" date serial from 45320 to 45392 ---
$wsAbout.Range("C1").Value = 45392

# --- MCF sheet: set capacity factor cells to 1 ---
$cellsToUpdate = @("B2", "B3", "B4", "B6", "B10", "B11", "B12", "B13", "B14", "B16", "B17", "B18")
foreach ($cellAddr in $cellsToUpdate) {
    $wsMCF.Range($cellAddr).Value = 1
}

# --- Update selection/active cell state on MCF sheet ---
$wsMCF.Activate()
$wsMCF.Range("B17").Select()

$wb.Save()
